$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture current (pre-edit) row values for rows 42-46 so we can
#     reshuffle them safely without clobbering data we still need ---
# (Column I is intentionally excluded: it is blank on every row in this
# sheet and must stay an untouched empty cell.)
function Get-RowData($row) {
    $data = @{}
    foreach ($col in @("A","B","D","E","F","G","H","M","P","Q","R","S","AC")) {
        $data[$col] = $ws.Range("$col$row").Value()
    }
    return $data
}

$row42 = Get-RowData 42
$row43 = Get-RowData 43
$row44 = Get-RowData 44
$row45 = Get-RowData 45
$row46 = Get-RowData 46

function Set-RowData($row, $data) {
    foreach ($col in @("A","B","D","E","F","G","H","M","P","Q","R","S","AC")) {
        $val = $data[$col]
        if ($val -eq $null) {
            $ws.Range("$col$row").Value = ""
        } else {
            $ws.Range("$col$row").Value = $val
        }
    }
}

# Row 42 becomes what row 44 used to be
Set-RowData 42 $row44

# Row 43 becomes what row 42 used to be, with the taxon sort order bumped
$row43New = $row42.Clone()
$row43New["B"] = 79244
Set-RowData 43 $row43New

# Row 44 becomes what row 43 used to be, with the taxon sort order bumped
$row44New = $row43.Clone()
$row44New["B"] = 79244
Set-RowData 44 $row44New

# Row 45 becomes what row 46 used to be
Set-RowData 45 $row46

# Row 46 becomes what row 45 used to be, with the taxon sort order bumped
$row46New = $row45.Clone()
$row46New["B"] = 79244
Set-RowData 46 $row46New

# --- Simple Taxonsorteringsordning bumps (79243 -> 79244) on the
#     remaining affected rows ---
$ws.Range("B39").Value = 79244
$ws.Range("B47").Value = 79244
$ws.Range("B49").Value = 79244
$ws.Range("B50").Value = 79244
$ws.Range("B52").Value = 79244
